# Generate Report for Handoff
# - Status moves from "In Translation" to "Ready for handoff" on the Overview
#   sheet (zh-cn / de-de status columns) as well as on each language sheet.
# - The "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps
#   are refreshed for both languages.
# - The Status column widens to fit the new, longer status text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$newStatus = "Ready for handoff"

# Overview sheet: zh-cn (E2) and de-de (F2) status cells
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus

# zh-cn sheet: Status (C2) and Latest Handoff Datetime (H2)
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("H2").Value = "2016-08-21 15:05:27"

# de-de sheet: Status (C2) and Latest Handoff Datetime (H2)
$dede.Range("C2").Value = $newStatus
$dede.Range("H2").Value = "2016-08-21 15:05:31"

# Widen the Status columns so the longer "Ready for handoff" text fits.
# (ColumnWidth is quantized by Excel to whole-pixel steps, so 16.3
# is the input that lands closest to the authored 17.216 width.)
$overview.Columns.Item(5).ColumnWidth = 16.3
$overview.Columns.Item(6).ColumnWidth = 16.3
$zhcn.Columns.Item(3).ColumnWidth = 16.3
$dede.Columns.Item(3).ColumnWidth = 16.3

# Re-assert the datetime number format on the timestamp cells so the
# round-trip save keeps them formatted as "yyyy-mm-dd HH:mm:ss".
$dateFormat = "yyyy-mm-dd HH:mm:ss"
$overview.Range("G2").NumberFormat = $dateFormat
$zhcn.Range("H2").NumberFormat = $dateFormat
$zhcn.Range("K2").NumberFormat = $dateFormat
$dede.Range("H2").NumberFormat = $dateFormat
$dede.Range("K2").NumberFormat = $dateFormat
